$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> ordered cell/value pairs reflecting the corrected financial figures
$rowData = @{}
$rowData[2] = @(@{Col="D"; Val=20207},@{Col="E"; Val=1623},@{Col="F"; Val=1623},@{Col="G"; Val=973},@{Col="H"; Val=1030},@{Col="I"; Val=1053},@{Col="J"; Val=-22},@{Col="K"; Val=30033},@{Col="L"; Val=17200},@{Col="M"; Val=12834},@{Col="N"; Val=12266},@{Col="O"; Val=568},@{Col="P"; Val=4015},@{Col="Q"; Val=2776},@{Col="R"; Val=-637},@{Col="S"; Val=-2136},@{Col="T"; Val=733},@{Col="U"; Val=2043},@{Col="V"; Val=11627},@{Col="W"; Val=8.029999999999999},@{Col="X"; Val=5.1},@{Col="Y"; Val=8.91},@{Col="Z"; Val=3.39},@{Col="AA"; Val=134.02},@{Col="AB"; Val=217.51},@{Col="AC"; Val=260},@{Col="AE"; Val=3034},@{Col="AF"; Val=0.79},@{Col="AG"; Val=0},@{Col="AH"; Val=0},@{Col="AI"; Val=0},@{Col="AJ"; Val=376140158})
$rowData[3] = @(@{Col="D"; Val=19864},@{Col="E"; Val=2142},@{Col="F"; Val=2206},@{Col="G"; Val=825},@{Col="H"; Val=771},@{Col="I"; Val=769},@{Col="J"; Val=2},@{Col="K"; Val=29345},@{Col="L"; Val=15862},@{Col="M"; Val=13483},@{Col="N"; Val=12926},@{Col="O"; Val=557},@{Col="P"; Val=4015},@{Col="Q"; Val=2742},@{Col="R"; Val=-127},@{Col="S"; Val=-2649},@{Col="T"; Val=648},@{Col="U"; Val=2093},@{Col="V"; Val=9473},@{Col="W"; Val=10.79},@{Col="X"; Val=3.88},@{Col="Y"; Val=6.1},@{Col="Z"; Val=2.6},@{Col="AA"; Val=117.65},@{Col="AB"; Val=233.57},@{Col="AC"; Val=190},@{Col="AE"; Val=3197},@{Col="AF"; Val=1.01},@{Col="AG"; Val=0},@{Col="AH"; Val=0},@{Col="AI"; Val=0},@{Col="AJ"; Val=376140158})
$rowData[4] = @(@{Col="D"; Val=14303},@{Col="E"; Val=2578},@{Col="F"; Val=2618},@{Col="G"; Val=2163},@{Col="H"; Val=1751},@{Col="I"; Val=1728},@{Col="J"; Val=23},@{Col="K"; Val=32162},@{Col="L"; Val=14022},@{Col="M"; Val=18140},@{Col="N"; Val=16889},@{Col="O"; Val=1251},@{Col="P"; Val=4656},@{Col="Q"; Val=1967},@{Col="R"; Val=-743},@{Col="S"; Val=511},@{Col="T"; Val=759},@{Col="U"; Val=1209},@{Col="V"; Val=7210},@{Col="W"; Val=18.03},@{Col="X"; Val=12.24},@{Col="Y"; Val=11.59},@{Col="Z"; Val=5.69},@{Col="AA"; Val=77.3},@{Col="AB"; Val=272.29},@{Col="AC"; Val=401},@{Col="AE"; Val=3628},@{Col="AF"; Val=0.82},@{Col="AG"; Val=32},@{Col="AH"; Val=1.08},@{Col="AI"; Val=16.21},@{Col="AJ"; Val=437394345})
$rowData[5] = @(@{Col="D"; Val=15171},@{Col="E"; Val=2509},@{Col="F"; Val=2509},@{Col="G"; Val=1760},@{Col="H"; Val=3012},@{Col="I"; Val=3021},@{Col="J"; Val=-8},@{Col="K"; Val=35012},@{Col="L"; Val=14801},@{Col="M"; Val=20211},@{Col="N"; Val=19604},@{Col="O"; Val=606},@{Col="P"; Val=5054},@{Col="Q"; Val=2983},@{Col="R"; Val=-2230},@{Col="S"; Val=47},@{Col="T"; Val=1364},@{Col="U"; Val=1619},@{Col="V"; Val=8799},@{Col="W"; Val=16.54},@{Col="X"; Val=19.86},@{Col="Y"; Val=16.55},@{Col="Z"; Val=8.970000000000001},@{Col="AA"; Val=73.23},@{Col="AB"; Val=289.06},@{Col="AC"; Val=613},@{Col="AD"; Val=6.08},@{Col="AE"; Val=3879},@{Col="AF"; Val=0.96},@{Col="AG"; Val=214},@{Col="AH"; Val=5.74},@{Col="AI"; Val=34.97},@{Col="AJ"; Val=503859595})
$rowData[6] = @(@{Col="D"; Val=15100},@{Col="E"; Val=2469},@{Col="F"; Val=2469},@{Col="G"; Val=2002},@{Col="H"; Val=1470},@{Col="I"; Val=1463},@{Col="K"; Val=34293},@{Col="L"; Val=14497},@{Col="M"; Val=19796},@{Col="N"; Val=19349},@{Col="P"; Val=5054},@{Col="Q"; Val=3160},@{Col="R"; Val=-1156},@{Col="S"; Val=-2569},@{Col="T"; Val=1355},@{Col="U"; Val=1805},@{Col="V"; Val=8969},@{Col="W"; Val=16.35},@{Col="X"; Val=9.74},@{Col="Y"; Val=7.51},@{Col="Z"; Val=4.24},@{Col="AA"; Val=73.23},@{Col="AB"; Val=285.85},@{Col="AC"; Val=290},@{Col="AD"; Val=21.73},@{Col="AE"; Val=3828},@{Col="AF"; Val=1.64},@{Col="AG"; Val=370},@{Col="AH"; Val=5.88},@{Col="AI"; Val=127.81},@{Col="AJ"; Val=503859595})
$rowData[7] = @(@{Col="D"; Val=15549},@{Col="E"; Val=2492},@{Col="G"; Val=2071},@{Col="H"; Val=1580},@{Col="I"; Val=1577},@{Col="K"; Val=34145},@{Col="L"; Val=14894},@{Col="M"; Val=19251},@{Col="N"; Val=18840},@{Col="P"; Val=5051},@{Col="Q"; Val=3441},@{Col="R"; Val=-1262},@{Col="S"; Val=-2040},@{Col="T"; Val=1219},@{Col="U"; Val=2518},@{Col="W"; Val=16.03},@{Col="X"; Val=10.16},@{Col="Y"; Val=8.26},@{Col="Z"; Val=4.62},@{Col="AA"; Val=77.36},@{Col="AC"; Val=312},@{Col="AD"; Val=15.92},@{Col="AE"; Val=3728},@{Col="AF"; Val=1.33},@{Col="AG"; Val=412},@{Col="AH"; Val=8.31},@{Col="AI"; Val=131.83})
$rowData[8] = @(@{Col="D"; Val=15715},@{Col="E"; Val=2583},@{Col="G"; Val=2207},@{Col="H"; Val=1666},@{Col="I"; Val=1661},@{Col="K"; Val=33626},@{Col="L"; Val=14882},@{Col="M"; Val=18745},@{Col="N"; Val=18356},@{Col="P"; Val=5051},@{Col="Q"; Val=3814},@{Col="R"; Val=-1403},@{Col="S"; Val=-2541},@{Col="T"; Val=1315},@{Col="U"; Val=2712},@{Col="W"; Val=16.44},@{Col="X"; Val=10.6},@{Col="Y"; Val=8.93},@{Col="Z"; Val=4.92},@{Col="AA"; Val=79.39},@{Col="AC"; Val=329},@{Col="AD"; Val=15.11},@{Col="AE"; Val=3632},@{Col="AF"; Val=1.37},@{Col="AG"; Val=429},@{Col="AH"; Val=8.640000000000001},@{Col="AI"; Val=130.04})
$rowData[9] = @(@{Col="D"; Val=16087},@{Col="E"; Val=2711},@{Col="G"; Val=2358},@{Col="H"; Val=1813},@{Col="I"; Val=1772},@{Col="K"; Val=33495},@{Col="L"; Val=15082},@{Col="M"; Val=18413},@{Col="N"; Val=17998},@{Col="P"; Val=5052},@{Col="Q"; Val=3943},@{Col="R"; Val=-1438},@{Col="S"; Val=-2487},@{Col="T"; Val=1320},@{Col="U"; Val=2845},@{Col="W"; Val=16.85},@{Col="X"; Val=11.27},@{Col="Y"; Val=9.75},@{Col="Z"; Val=5.4},@{Col="AA"; Val=81.91},@{Col="AC"; Val=351},@{Col="AD"; Val=14.16},@{Col="AE"; Val=3561},@{Col="AF"; Val=1.39},@{Col="AG"; Val=429},@{Col="AH"; Val=8.630000000000001},@{Col="AI"; Val=121.83})

foreach ($r in $rowData.Keys) {
    foreach ($item in $rowData[$r]) {
        $ws.Range("$($item.Col)$r").Value = $item.Val
    }
}

# Column AD is dropped entirely for rows 2-4 in the corrected data
$ws.Range("AD2").ClearContents()
$ws.Range("AD3").ClearContents()
$ws.Range("AD4").ClearContents()
